# Apply "2022-02-22 data" update to the Fonds de solidarite (volet 1)
# regional / classe effectif workbook: updates nombre_aides (col C) and
# montant_total (col E) for a set of rows to their refreshed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Row = 2;   C = 100800; E = 327272510 }
    @{ Row = 3;   C = 249273; E = 1036172372 }
    @{ Row = 5;   C = 39477;  E = 361411392 }
    @{ Row = 46;  C = 10920;  E = 66193511 }
    @{ Row = 53;  C = 141655; E = 589981563 }
    @{ Row = 55;  C = 23187;  E = 187773972 }
    @{ Row = 57;  C = 3705;   E = 137905354 }
    @{ Row = 63;  C = 14101;  E = 35459928 }
    @{ Row = 64;  C = 5055;   E = 19389984 }
    @{ Row = 79;  C = 116574; E = 447276129 }
    @{ Row = 81;  C = 17423;  E = 133403196 }
    @{ Row = 91;  C = 150984; E = 480931655 }
    @{ Row = 92;  C = 408705; E = 1590402401 }
    @{ Row = 93;  C = 209256; E = 1303951700 }
    @{ Row = 94;  C = 93999;  E = 911298136 }
    @{ Row = 96;  C = 17144;  E = 780902892 }
    @{ Row = 104; C = 135152; E = 271713428 }
    @{ Row = 106; C = 18118;  E = 40712935 }
    @{ Row = 114; C = 3714;   E = 8915019 }
    @{ Row = 115; C = 11468;  E = 32283158 }
    @{ Row = 116; C = 4424;   E = 19530976 }
    @{ Row = 118; C = 906;    E = 10421491 }
    @{ Row = 122; C = 8324;   E = 12555568 }
    @{ Row = 131; C = 75573;  E = 307076240 }
    @{ Row = 166; C = 35925;  E = 210542620 }
    @{ Row = 174; C = 226046; E = 900132284 }
    @{ Row = 175; C = 80749;  E = 485283850 }
    @{ Row = 184; C = 68726;  E = 134113858 }
)

foreach ($change in $changes) {
    $ws.Cells.Item($change.Row, 3).Value = $change.C
    $ws.Cells.Item($change.Row, 5).Value = $change.E
}
